$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"

$ws.Range("D2").Value = '54.318.19'
$ws.Range("E2").Value = '  +1.34%  '

$ws.Range("D3").Value = '2.272.25'
$ws.Range("E3").Value = '  +1.46%  '

$ws.Range("E4").Value = '  +0.21%  '

$ws.Range("D5").Value = '496.53'
$ws.Range("E5").Value = '  +1.96%  '

$ws.Range("D6").Value = '128.20'
$ws.Range("E6").Value = '  +2.29%  '

$ws.Range("E7").Value = '  +0.19%  '

$ws.Range("D8").Value = '0.528'
$ws.Range("E8").Value = '  +1.06%  '

$ws.Range("D9").Value = '0.0958'
$ws.Range("E9").Value = '  +3.98%  '

$ws.Range("E10").Value = '  +2.12%  '

$ws.Range("D11").Value = '0.331'
$ws.Range("E11").Value = '  +3.62%  '

$ws.Range("D12").Value = '4.70'
$ws.Range("E12").Value = '  +1.64%  '

$ws.Range("D13").Value = '2.675.65'
$ws.Range("E13").Value = '  +2.48%  '

$ws.Range("D14").Value = '22.36'
$ws.Range("E14").Value = '  +5.28%  '

$ws.Range("D15").Value = '54.245.46'
$ws.Range("E15").Value = '  +1.32%  '

$ws.Range("E16").Value = '  +1.44%  '

$ws.Range("D17").Value = '2.276.80'
$ws.Range("E17").Value = '  +2.42%  '

$ws.Range("D18").Value = '10.14'
$ws.Range("E18").Value = '  +5.04%  '

$ws.Range("E19").Value = '  +3.39%  '

$ws.Range("D20").Value = '304.45'
$ws.Range("E20").Value = '  +2.99%  '

$ws.Range("E21").Value = '  +4.74%  '

$ws.Range("E22").Value = '  +0.50%  '

$ws.Range("D23").Value = '61.90'
$ws.Range("E23").Value = '  -2.74%  '

$ws.Range("E24").Value = '  +0.12%  '

$ws.Range("D25").Value = '2.377.65'
$ws.Range("E25").Value = '  +2.68%  '

$ws.Range("E26").Value = '  +2.18%  '

$ws.Range("D27").Value = '7.25'
$ws.Range("E27").Value = '  +2.96%  '

$ws.Range("D28").Value = '171.59'
$ws.Range("E28").Value = '  +4.77%  '

$ws.Range("E29").Value = '  +2.16%  '

$ws.Range("D30").Value = '0.0₃0683'
$ws.Range("E30").Value = '  +2.12%  '

$ws.Range("D31").Value = '5.90'
$ws.Range("E31").Value = '  +2.16%  '

$ws.Range("E32").Value = '  +2.92%  '

$ws.Range("E33").Value = '  +0.14%  '

$ws.Range("D34").Value = '17.74'
$ws.Range("E34").Value = '  +2.71%  '

$ws.Range("E35").Value = '  -0.03%  '

$ws.Range("D36").Value = '0.919'
$ws.Range("E36").Value = '  +10.47%  '

$ws.Range("D37").Value = '1.20'
$ws.Range("E37").Value = '  +2.48%  '

$ws.Range("E38").Value = '  +3.99%  '

$ws.Range("D39").Value = '35.75'
$ws.Range("E39").Value = '  +1.63%  '

$ws.Range("D40").Value = '0.373'
$ws.Range("E40").Value = '  +1.18%  '

$ws.Range("E41").Value = '  +2.50%  '

$ws.Range("E42").Value = '  +3.03%  '

$ws.Range("B43").Value = 'RenderToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D43").Value = '4.98'
$ws.Range("E43").Value = '  +3.35%  '

$ws.Range("B44").Value = 'Aave'
$ws.Range("C44").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D44").Value = '125.87'
$ws.Range("E44").Value = '  -0.50%  '

$ws.Range("E45").Value = '  +1.93%  '

$ws.Range("D46").Value = '0.0490'
$ws.Range("E46").Value = '  +3.97%  '

$ws.Range("D47").Value = '0.546'
$ws.Range("E47").Value = '  +2.07%  '

$ws.Range("D48").Value = '237.81'
$ws.Range("E48").Value = '  +1.37%  '

$ws.Range("E49").Value = '  +1.03%  '

$ws.Range("D50").Value = '0.0206'
$ws.Range("E50").Value = '  +2.63%  '

$ws.Range("E51").Value = '  +0.97%  '

$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").Style = "Normal"
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").Style = "Normal"
$ws.Range("D14").Style = "Normal"
$ws.Range("D18").Style = "Normal"
$ws.Range("D20").Style = "Normal"
$ws.Range("D23").Style = "Normal"
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").Style = "Normal"
$ws.Range("D31").Style = "Normal"
$ws.Range("D34").Style = "Normal"
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").Style = "Normal"
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").Style = "Normal"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").Style = "Normal"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").Style = "Normal"
$ws.Range("D50").Style = "Normal"